# Calcs and values updated for the new AFE approach
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, pushing old rows 5-9 down to 6-10.
# Excel will automatically adjust the relative formulas in the shifted rows.
$ws.Rows("5:5").Insert()

# New row 6: "Amp. Op. Gain" label + gain values used by the new AFE formulas
$ws.Range("B6").Value = "Amp. Op. Gain"
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 2

# Column B now needs to be wide enough to show the new, longer label
$ws.Columns("B:B").AutoFit()

# Update the reference voltage (C4) for the new AFE approach, and make E4
# simply mirror C4 via a formula instead of a hard-coded duplicate value.
$ws.Range("C4").Value = 1.6
$ws.Range("E4").Formula = '=C4'

# Update the shifted Vout formulas (now on row 7) to account for the new gain
$ws.Range("C7").Formula = '=(($C$4/($C$3+$C$4))*C2)*2'
$ws.Range("E7").Formula = '=(($C$4/($C$3+$C$4))*E2)*2'

# Update the window view settings recorded for the workbook
$win = $excel.ActiveWindow
$win.Width = 19200
$win.Height = 10665
$win.Left = -9165
$win.Top = 3195

# Update the active selection to match the new authored state
$ws.Range("C5").Select()

$wb.Save()
